$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Cells.Item(4, 2).Value = 0.296
$ws.Cells.Item(4, 3).Value = 0.051
$ws.Cells.Item(4, 5).Value = 0.151
$ws.Cells.Item(4, 8).Value = 0.189
$ws.Cells.Item(4, 10).Value = 0.108
$ws.Cells.Item(4, 11).Value = 0.342
$ws.Cells.Item(4, 12).Value = 0.102
$ws.Cells.Item(4, 13).Value = 0.32
$ws.Cells.Item(4, 14).Value = 0.271
$ws.Cells.Item(4, 15).Value = 0.02
$ws.Cells.Item(4, 16).Value = 0.143
$ws.Cells.Item(4, 17).Value = 0.513
$ws.Cells.Item(4, 18).Value = 0.217
$ws.Cells.Item(4, 19).Value = 0.466
$ws.Cells.Item(4, 20).Value = 0.283
$ws.Cells.Item(4, 23).Value = 0.244
$ws.Cells.Item(4, 25).Value = 0.208
$ws.Cells.Item(4, 26).Value = 0.451
$ws.Cells.Item(4, 27).Value = 0.133
$ws.Cells.Item(4, 28).Value = 0.364
$ws.Cells.Item(4, 29).Value = 0.126
$ws.Cells.Item(4, 31).Value = 0.078
$ws.Cells.Item(4, 32).Value = 0.737
$ws.Cells.Item(4, 33).Value = 0.094
$ws.Cells.Item(4, 34).Value = 0.307
$ws.Cells.Item(4, 35).Value = 0.658
$ws.Cells.Item(4, 36).Value = 0.172
$ws.Cells.Item(4, 37).Value = 0.415
$ws.Cells.Item(4, 38).Value = 0.703
$ws.Cells.Item(4, 40).Value = 0.339
$ws.Cells.Item(4, 41).Value = 0.699

# Row 5
$ws.Cells.Item(5, 2).Value = 0.8159999999999999
$ws.Cells.Item(5, 3).Value = 0.15
$ws.Cells.Item(5, 4).Value = 0.388
$ws.Cells.Item(5, 5).Value = 0.658
$ws.Cells.Item(5, 6).Value = 0.225
$ws.Cells.Item(5, 7).Value = 0.474
$ws.Cells.Item(5, 8).Value = 0.8159999999999999
$ws.Cells.Item(5, 9).Value = 0.15
$ws.Cells.Item(5, 10).Value = 0.388
$ws.Cells.Item(5, 11).Value = 0.658
$ws.Cells.Item(5, 12).Value = 0.225
$ws.Cells.Item(5, 13).Value = 0.474
$ws.Cells.Item(5, 14).Value = 0.842
$ws.Cells.Item(5, 15).Value = 0.133
$ws.Cells.Item(5, 16).Value = 0.365
$ws.Cells.Item(5, 17).Value = 0.579
$ws.Cells.Item(5, 18).Value = 0.244
$ws.Cells.Item(5, 19).Value = 0.494
$ws.Cells.Item(5, 20).Value = 0.579
$ws.Cells.Item(5, 21).Value = 0.244
$ws.Cells.Item(5, 22).Value = 0.494
$ws.Cells.Item(5, 23).Value = 0.737
$ws.Cells.Item(5, 24).Value = 0.194
$ws.Cells.Item(5, 25).Value = 0.44
$ws.Cells.Item(5, 26).Value = 0.8159999999999999
$ws.Cells.Item(5, 27).Value = 0.15
$ws.Cells.Item(5, 28).Value = 0.388
$ws.Cells.Item(5, 29).Value = 0.763
$ws.Cells.Item(5, 30).Value = 0.181
$ws.Cells.Item(5, 31).Value = 0.425
$ws.Cells.Item(5, 32).Value = 0.974
$ws.Cells.Item(5, 34).Value = 0.16
$ws.Cells.Item(5, 35).Value = 0.763
$ws.Cells.Item(5, 36).Value = 0.181
$ws.Cells.Item(5, 37).Value = 0.425
$ws.Cells.Item(5, 38).Value = 0.921
$ws.Cells.Item(5, 39).Value = 0.073
$ws.Cells.Item(5, 40).Value = 0.27
$ws.Cells.Item(5, 41).Value = 0.886

# Row 6
$ws.Cells.Item(6, 2).Value = 0.434
$ws.Cells.Item(6, 5).Value = 0.246
$ws.Cells.Item(6, 8).Value = 0.307
$ws.Cells.Item(6, 11).Value = 0.45
$ws.Cells.Item(6, 14).Value = 0.41
$ws.Cells.Item(6, 17).Value = 0.544
$ws.Cells.Item(6, 20).Value = 0.38
$ws.Cells.Item(6, 23).Value = 0.367
$ws.Cells.Item(6, 26).Value = 0.581
$ws.Cells.Item(6, 29).Value = 0.216
$ws.Cells.Item(6, 32).Value = 0.839
$ws.Cells.Item(6, 35).Value = 0.707
$ws.Cells.Item(6, 38).Value = 0.797
$ws.Cells.Item(6, 41).Value = 0.781

# Row 7
$ws.Cells.Item(7, 2).Value = 0.604
$ws.Cells.Item(7, 5).Value = 0.394
$ws.Cells.Item(7, 8).Value = 0.491
$ws.Cells.Item(7, 11).Value = 0.555
$ws.Cells.Item(7, 14).Value = 0.592
$ws.Cells.Item(7, 17).Value = 0.5639999999999999
$ws.Cells.Item(7, 20).Value = 0.479
$ws.Cells.Item(7, 23).Value = 0.525
$ws.Cells.Item(7, 26).Value = 0.702
$ws.Cells.Item(7, 29).Value = 0.379
$ws.Cells.Item(7, 32).Value = 0.915
$ws.Cells.Item(7, 35).Value = 0.739
$ws.Cells.Item(7, 38).Value = 0.867
$ws.Cells.Item(7, 41).Value = 0.84

# Row 8
$ws.Cells.Item(8, 2).Value = 0.751
$ws.Cells.Item(8, 3).Value = 0.15
$ws.Cells.Item(8, 4).Value = 0.387
$ws.Cells.Item(8, 5).Value = 0.548
$ws.Cells.Item(8, 6).Value = 0.191
$ws.Cells.Item(8, 7).Value = 0.437
$ws.Cells.Item(8, 8).Value = 0.703
$ws.Cells.Item(8, 9).Value = 0.152
$ws.Cells.Item(8, 10).Value = 0.39
$ws.Cells.Item(8, 11).Value = 0.586
$ws.Cells.Item(8, 12).Value = 0.204
$ws.Cells.Item(8, 13).Value = 0.452
$ws.Cells.Item(8, 14).Value = 0.749
$ws.Cells.Item(8, 15).Value = 0.137
$ws.Cells.Item(8, 16).Value = 0.371
$ws.Cells.Item(8, 17).Value = 0.55
$ws.Cells.Item(8, 18).Value = 0.229
$ws.Cells.Item(8, 19).Value = 0.479
$ws.Cells.Item(8, 20).Value = 0.501
$ws.Cells.Item(8, 21).Value = 0.206
$ws.Cells.Item(8, 22).Value = 0.454
$ws.Cells.Item(8, 23).Value = 0.662
$ws.Cells.Item(8, 24).Value = 0.18
$ws.Cells.Item(8, 25).Value = 0.424
$ws.Cells.Item(8, 26).Value = 0.751
$ws.Cells.Item(8, 27).Value = 0.15
$ws.Cells.Item(8, 28).Value = 0.387
$ws.Cells.Item(8, 29).Value = 0.655
$ws.Cells.Item(8, 30).Value = 0.177
$ws.Cells.Item(8, 31).Value = 0.42
$ws.Cells.Item(8, 32).Value = 0.893
$ws.Cells.Item(8, 33).Value = 0.046
$ws.Cells.Item(8, 34).Value = 0.215
$ws.Cells.Item(8, 35).Value = 0.753
$ws.Cells.Item(8, 36).Value = 0.18
$ws.Cells.Item(8, 37).Value = 0.424
$ws.Cells.Item(8, 38).Value = 0.892
$ws.Cells.Item(8, 39).Value = 0.078
$ws.Cells.Item(8, 40).Value = 0.279
$ws.Cells.Item(8, 41).Value = 0.846

# Row 9
$ws.Cells.Item(9, 2).Value = 0.658
$ws.Cells.Item(9, 3).Value = 0.225
$ws.Cells.Item(9, 4).Value = 0.474
$ws.Cells.Item(9, 5).Value = 0.421
$ws.Cells.Item(9, 6).Value = 0.244
$ws.Cells.Item(9, 7).Value = 0.494
$ws.Cells.Item(9, 8).Value = 0.579
$ws.Cells.Item(9, 9).Value = 0.244
$ws.Cells.Item(9, 10).Value = 0.494
$ws.Cells.Item(9, 11).Value = 0.5
$ws.Cells.Item(9, 14).Value = 0.632
$ws.Cells.Item(9, 15).Value = 0.233
$ws.Cells.Item(9, 16).Value = 0.482
$ws.Cells.Item(9, 17).Value = 0.5
$ws.Cells.Item(9, 20).Value = 0.395
$ws.Cells.Item(9, 21).Value = 0.239
$ws.Cells.Item(9, 22).Value = 0.489
$ws.Cells.Item(9, 23).Value = 0.553
$ws.Cells.Item(9, 24).Value = 0.247
$ws.Cells.Item(9, 25).Value = 0.497
$ws.Cells.Item(9, 26).Value = 0.658
$ws.Cells.Item(9, 27).Value = 0.225
$ws.Cells.Item(9, 28).Value = 0.474
$ws.Cells.Item(9, 29).Value = 0.553
$ws.Cells.Item(9, 30).Value = 0.247
$ws.Cells.Item(9, 31).Value = 0.497
$ws.Cells.Item(9, 32).Value = 0.763
$ws.Cells.Item(9, 33).Value = 0.181
$ws.Cells.Item(9, 34).Value = 0.425
$ws.Cells.Item(9, 35).Value = 0.737
$ws.Cells.Item(9, 36).Value = 0.194
$ws.Cells.Item(9, 37).Value = 0.44
$ws.Cells.Item(9, 38).Value = 0.842
$ws.Cells.Item(9, 39).Value = 0.133
$ws.Cells.Item(9, 40).Value = 0.365
$ws.Cells.Item(9, 41).Value = 0.781

# Row 10
$ws.Cells.Item(10, 2).Value = 0.8159999999999999
$ws.Cells.Item(10, 3).Value = 0.15
$ws.Cells.Item(10, 4).Value = 0.388
$ws.Cells.Item(10, 5).Value = 0.579
$ws.Cells.Item(10, 6).Value = 0.244
$ws.Cells.Item(10, 7).Value = 0.494
$ws.Cells.Item(10, 8).Value = 0.737
$ws.Cells.Item(10, 9).Value = 0.194
$ws.Cells.Item(10, 10).Value = 0.44
$ws.Cells.Item(10, 11).Value = 0.658
$ws.Cells.Item(10, 12).Value = 0.225
$ws.Cells.Item(10, 13).Value = 0.474
$ws.Cells.Item(10, 14).Value = 0.8159999999999999
$ws.Cells.Item(10, 15).Value = 0.15
$ws.Cells.Item(10, 16).Value = 0.388
$ws.Cells.Item(10, 17).Value = 0.579
$ws.Cells.Item(10, 18).Value = 0.244
$ws.Cells.Item(10, 19).Value = 0.494
$ws.Cells.Item(10, 20).Value = 0.579
$ws.Cells.Item(10, 21).Value = 0.244
$ws.Cells.Item(10, 22).Value = 0.494
$ws.Cells.Item(10, 23).Value = 0.737
$ws.Cells.Item(10, 24).Value = 0.194
$ws.Cells.Item(10, 25).Value = 0.44
$ws.Cells.Item(10, 26).Value = 0.8159999999999999
$ws.Cells.Item(10, 27).Value = 0.15
$ws.Cells.Item(10, 28).Value = 0.388
$ws.Cells.Item(10, 29).Value = 0.658
$ws.Cells.Item(10, 30).Value = 0.225
$ws.Cells.Item(10, 31).Value = 0.474
$ws.Cells.Item(10, 32).Value = 0.974
$ws.Cells.Item(10, 34).Value = 0.16
$ws.Cells.Item(10, 35).Value = 0.763
$ws.Cells.Item(10, 36).Value = 0.181
$ws.Cells.Item(10, 37).Value = 0.425
$ws.Cells.Item(10, 38).Value = 0.921
$ws.Cells.Item(10, 39).Value = 0.073
$ws.Cells.Item(10, 40).Value = 0.27
$ws.Cells.Item(10, 41).Value = 0.886

# Row 11
$ws.Cells.Item(11, 2).Value = 0.8159999999999999
$ws.Cells.Item(11, 3).Value = 0.15
$ws.Cells.Item(11, 4).Value = 0.388
$ws.Cells.Item(11, 5).Value = 0.658
$ws.Cells.Item(11, 6).Value = 0.225
$ws.Cells.Item(11, 7).Value = 0.474
$ws.Cells.Item(11, 8).Value = 0.8159999999999999
$ws.Cells.Item(11, 9).Value = 0.15
$ws.Cells.Item(11, 10).Value = 0.388
$ws.Cells.Item(11, 11).Value = 0.658
$ws.Cells.Item(11, 12).Value = 0.225
$ws.Cells.Item(11, 13).Value = 0.474
$ws.Cells.Item(11, 14).Value = 0.842
$ws.Cells.Item(11, 15).Value = 0.133
$ws.Cells.Item(11, 16).Value = 0.365
$ws.Cells.Item(11, 17).Value = 0.579
$ws.Cells.Item(11, 18).Value = 0.244
$ws.Cells.Item(11, 19).Value = 0.494
$ws.Cells.Item(11, 20).Value = 0.579
$ws.Cells.Item(11, 21).Value = 0.244
$ws.Cells.Item(11, 22).Value = 0.494
$ws.Cells.Item(11, 23).Value = 0.737
$ws.Cells.Item(11, 24).Value = 0.194
$ws.Cells.Item(11, 25).Value = 0.44
$ws.Cells.Item(11, 26).Value = 0.8159999999999999
$ws.Cells.Item(11, 27).Value = 0.15
$ws.Cells.Item(11, 28).Value = 0.388
$ws.Cells.Item(11, 29).Value = 0.711
$ws.Cells.Item(11, 30).Value = 0.206
$ws.Cells.Item(11, 31).Value = 0.454
$ws.Cells.Item(11, 32).Value = 0.974
$ws.Cells.Item(11, 34).Value = 0.16
$ws.Cells.Item(11, 35).Value = 0.763
$ws.Cells.Item(11, 36).Value = 0.181
$ws.Cells.Item(11, 37).Value = 0.425
$ws.Cells.Item(11, 38).Value = 0.921
$ws.Cells.Item(11, 39).Value = 0.073
$ws.Cells.Item(11, 40).Value = 0.27
$ws.Cells.Item(11, 41).Value = 0.886

# Row 12
$ws.Cells.Item(12, 2).Value = 1.258
$ws.Cells.Item(12, 3).Value = 0.32
$ws.Cells.Item(12, 4).Value = 0.5659999999999999
$ws.Cells.Item(12, 5).Value = 1.68
$ws.Cells.Item(12, 6).Value = 1.098
$ws.Cells.Item(12, 7).Value = 1.048
$ws.Cells.Item(12, 8).Value = 1.613
$ws.Cells.Item(12, 9).Value = 1.334
$ws.Cells.Item(12, 10).Value = 1.155
$ws.Cells.Item(12, 11).Value = 1.4
$ws.Cells.Item(12, 12).Value = 0.5600000000000001
$ws.Cells.Item(12, 13).Value = 0.748
$ws.Cells.Item(12, 14).Value = 1.406
$ws.Cells.Item(12, 15).Value = 0.616
$ws.Cells.Item(12, 16).Value = 0.785
$ws.Cells.Item(12, 26).Value = 1.258
$ws.Cells.Item(12, 27).Value = 0.32
$ws.Cells.Item(12, 28).Value = 0.5659999999999999
$ws.Cells.Item(12, 29).Value = 1.793
$ws.Cells.Item(12, 30).Value = 2.44
$ws.Cells.Item(12, 31).Value = 1.562
$ws.Cells.Item(12, 32).Value = 1.243
$ws.Cells.Item(12, 33).Value = 0.238
$ws.Cells.Item(12, 34).Value = 0.488
$ws.Cells.Item(12, 35).Value = 1.034
$ws.Cells.Item(12, 36).Value = 0.033
$ws.Cells.Item(12, 37).Value = 0.182
$ws.Cells.Item(12, 38).Value = 1.086
$ws.Cells.Item(12, 39).Value = 0.078
$ws.Cells.Item(12, 40).Value = 0.28
$ws.Cells.Item(12, 41).Value = 1.121

# Row 13
$ws.Cells.Item(13, 2).Value = 3.474
$ws.Cells.Item(13, 3).Value = 1.46
$ws.Cells.Item(13, 4).Value = 1.208
$ws.Cells.Item(13, 5).Value = 4.594
$ws.Cells.Item(13, 6).Value = 0.429
$ws.Cells.Item(13, 7).Value = 0.655
$ws.Cells.Item(13, 8).Value = 4.611
$ws.Cells.Item(13, 9).Value = 0.627
$ws.Cells.Item(13, 10).Value = 0.792
$ws.Cells.Item(13, 11).Value = 2.265
$ws.Cells.Item(13, 12).Value = 0.606
$ws.Cells.Item(13, 13).Value = 0.779
$ws.Cells.Item(13, 14).Value = 3.263
$ws.Cells.Item(13, 15).Value = 0.72
$ws.Cells.Item(13, 16).Value = 0.849
$ws.Cells.Item(13, 26).Value = 2.514
$ws.Cells.Item(13, 27).Value = 2.878
$ws.Cells.Item(13, 28).Value = 1.697
$ws.Cells.Item(13, 29).Value = 6.378
$ws.Cells.Item(13, 30).Value = 2.181
$ws.Cells.Item(13, 31).Value = 1.477
$ws.Cells.Item(13, 32).Value = 1.605
$ws.Cells.Item(13, 33).Value = 0.713
$ws.Cells.Item(13, 34).Value = 0.844
$ws.Cells.Item(13, 35).Value = 1.289
$ws.Cells.Item(13, 36).Value = 0.364
$ws.Cells.Item(13, 37).Value = 0.603
$ws.Cells.Item(13, 38).Value = 1.579
$ws.Cells.Item(13, 39).Value = 0.717
$ws.Cells.Item(13, 40).Value = 0.847
$ws.Cells.Item(13, 41).Value = 1.491

